# Scheduled-runner style update of market price / profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
#  LeveProfitNQ/HQ -> columns H..N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets, reflecting refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 126.30769
$ws.Cells.Item(9, 9).Value = 150.625
$ws.Cells.Item(9, 10).Value = 87.40000000000001
$ws.Cells.Item(9, 11).Value = 150.625
$ws.Cells.Item(9, 12).Value = 87.40000000000001
$ws.Cells.Item(9, 13).Value = 18.375
$ws.Cells.Item(9, 14).Value = -425.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 624.4483
$ws.Cells.Item(18, 9).Value = 648.4815
$ws.Cells.Item(18, 11).Value = 648.4815
$ws.Cells.Item(18, 13).Value = -364.4815

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 1513.4615
$ws.Cells.Item(33, 9).Value = 148.75
$ws.Cells.Item(33, 10).Value = 3697
$ws.Cells.Item(33, 11).Value = 148.75
$ws.Cells.Item(33, 12).Value = 3697
$ws.Cells.Item(33, 13).Value = 80.25
$ws.Cells.Item(33, 14).Value = -4155

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 854.61536
$ws.Cells.Item(96, 9).Value = 619.2857
$ws.Cells.Item(96, 10).Value = 1129.1666
$ws.Cells.Item(96, 11).Value = 1857.8571
$ws.Cells.Item(96, 12).Value = 3387.4998
$ws.Cells.Item(96, 13).Value = -484.8571000000002
$ws.Cells.Item(96, 14).Value = -6133.4998

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 904
$ws.Cells.Item(98, 9).Value = 869.6667
$ws.Cells.Item(98, 10).Value = 972.6667
$ws.Cells.Item(98, 11).Value = 869.6667
$ws.Cells.Item(98, 12).Value = 972.6667
$ws.Cells.Item(98, 13).Value = 628.3333
$ws.Cells.Item(98, 14).Value = -3968.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 904
$ws.Cells.Item(122, 9).Value = 869.6667
$ws.Cells.Item(122, 10).Value = 972.6667
$ws.Cells.Item(122, 11).Value = 2609.0001
$ws.Cells.Item(122, 12).Value = 2918.0001
$ws.Cells.Item(122, 13).Value = -159.0001000000002
$ws.Cells.Item(122, 14).Value = -7818.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1596
$ws.Cells.Item(137, 9).Value = 1100
$ws.Cells.Item(137, 10).Value = 1676
$ws.Cells.Item(137, 11).Value = 3300
$ws.Cells.Item(137, 12).Value = 5028
$ws.Cells.Item(137, 13).Value = -750
$ws.Cells.Item(137, 14).Value = -10128

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 36056.74
$ws.Cells.Item(32, 9).Value = 6949.8823
$ws.Cells.Item(32, 10).Value = 97908.81
$ws.Cells.Item(32, 11).Value = 6949.8823
$ws.Cells.Item(32, 12).Value = 97908.81
$ws.Cells.Item(32, 13).Value = -6662.8823
$ws.Cells.Item(32, 14).Value = -98482.81

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1418.45
$ws.Cells.Item(45, 9).Value = 1432.4445
$ws.Cells.Item(45, 10).Value = 1407
$ws.Cells.Item(45, 11).Value = 1432.4445
$ws.Cells.Item(45, 12).Value = 1407
$ws.Cells.Item(45, 13).Value = -1055.4445
$ws.Cells.Item(45, 14).Value = -2161

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 1486
$ws.Cells.Item(80, 9).Value = 623.8570999999999
$ws.Cells.Item(80, 11).Value = 623.8570999999999
$ws.Cells.Item(80, 13).Value = 374.1429000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 1486
$ws.Cells.Item(83, 9).Value = 623.8570999999999
$ws.Cells.Item(83, 11).Value = 3119.2855
$ws.Cells.Item(83, 13).Value = 1872.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 979.5
$ws.Cells.Item(16, 9).Value = 788.2
$ws.Cells.Item(16, 10).Value = 1298.3334
$ws.Cells.Item(16, 11).Value = 788.2
$ws.Cells.Item(16, 12).Value = 1298.3334
$ws.Cells.Item(16, 13).Value = -501.2
$ws.Cells.Item(16, 14).Value = -1872.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 9938.5
$ws.Cells.Item(41, 9).Value = 59
$ws.Cells.Item(41, 10).Value = 10698.462
$ws.Cells.Item(41, 11).Value = 59
$ws.Cells.Item(41, 12).Value = 10698.462
$ws.Cells.Item(41, 13).Value = 369
$ws.Cells.Item(41, 14).Value = -11554.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 10730
$ws.Cells.Item(50, 10).Value = 10730
$ws.Cells.Item(50, 12).Value = 10730
$ws.Cells.Item(50, 14).Value = -11980

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 7173.8
$ws.Cells.Item(51, 10).Value = 7949.75
$ws.Cells.Item(51, 12).Value = 7949.75
$ws.Cells.Item(51, 14).Value = -9421.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 24846.666
$ws.Cells.Item(59, 10).Value = 24846.666
$ws.Cells.Item(59, 12).Value = 24846.666
$ws.Cells.Item(59, 14).Value = -27136.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(60, 8).Value = 17470
$ws.Cells.Item(60, 10).Value = 19293.334
$ws.Cells.Item(60, 12).Value = 19293.334
$ws.Cells.Item(60, 14).Value = -20315.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(61, 8).Value = 7173.8
$ws.Cells.Item(61, 10).Value = 7949.75
$ws.Cells.Item(61, 12).Value = 7949.75
$ws.Cells.Item(61, 14).Value = -8645.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3825.818
$ws.Cells.Item(62, 9).Value = 2501
$ws.Cells.Item(62, 10).Value = 4582.857
$ws.Cells.Item(62, 11).Value = 2501
$ws.Cells.Item(62, 12).Value = 4582.857
$ws.Cells.Item(62, 13).Value = -1877
$ws.Cells.Item(62, 14).Value = -5830.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 3825.818
$ws.Cells.Item(65, 9).Value = 2501
$ws.Cells.Item(65, 10).Value = 4582.857
$ws.Cells.Item(65, 11).Value = 12505
$ws.Cells.Item(65, 12).Value = 22914.285
$ws.Cells.Item(65, 13).Value = -9385
$ws.Cells.Item(65, 14).Value = -29154.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 888.4666999999999
$ws.Cells.Item(94, 9).Value = 850
$ws.Cells.Item(94, 10).Value = 891.2143
$ws.Cells.Item(94, 11).Value = 850
$ws.Cells.Item(94, 12).Value = 891.2143
$ws.Cells.Item(94, 13).Value = -399
$ws.Cells.Item(94, 14).Value = -1793.2143

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 996.7619
$ws.Cells.Item(105, 9).Value = 940.6667
$ws.Cells.Item(105, 10).Value = 1333.3334
$ws.Cells.Item(105, 11).Value = 940.6667
$ws.Cells.Item(105, 12).Value = 1333.3334
$ws.Cells.Item(105, 13).Value = 806.3333
$ws.Cells.Item(105, 14).Value = -4827.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 776.6
$ws.Cells.Item(107, 9).Value = 943.44446
$ws.Cells.Item(107, 10).Value = 640.0909
$ws.Cells.Item(107, 11).Value = 943.44446
$ws.Cells.Item(107, 12).Value = 640.0909
$ws.Cells.Item(107, 13).Value = 976.55554
$ws.Cells.Item(107, 14).Value = -4480.0909

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 979.5
$ws.Cells.Item(113, 9).Value = 788.2
$ws.Cells.Item(113, 10).Value = 1298.3334
$ws.Cells.Item(113, 11).Value = 788.2
$ws.Cells.Item(113, 12).Value = 1298.3334
$ws.Cells.Item(113, 13).Value = 1381.8
$ws.Cells.Item(113, 14).Value = -5638.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2380.16
$ws.Cells.Item(122, 9).Value = 2720.2856
$ws.Cells.Item(122, 10).Value = 1947.2727
$ws.Cells.Item(122, 11).Value = 8160.8568
$ws.Cells.Item(122, 12).Value = 5841.8181
$ws.Cells.Item(122, 13).Value = -5710.8568
$ws.Cells.Item(122, 14).Value = -10741.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1836.7858
$ws.Cells.Item(134, 9).Value = 1414.7273
$ws.Cells.Item(134, 11).Value = 4244.1819
$ws.Cells.Item(134, 13).Value = -1709.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 39.789474
$ws.Cells.Item(12, 10).Value = 47.533333
$ws.Cells.Item(12, 12).Value = 142.599999
$ws.Cells.Item(12, 14).Value = -488.599999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2540.7896
$ws.Cells.Item(131, 10).Value = 2767.353
$ws.Cells.Item(131, 12).Value = 8302.059000000001
$ws.Cells.Item(131, 14).Value = -18382.059

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 81392.66
$ws.Cells.Item(70, 9).Value = 129310.125
$ws.Cells.Item(70, 10).Value = 4724.7
$ws.Cells.Item(70, 11).Value = 129310.125
$ws.Cells.Item(70, 12).Value = 4724.7
$ws.Cells.Item(70, 13).Value = -129040.125
$ws.Cells.Item(70, 14).Value = -5264.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 81392.66
$ws.Cells.Item(73, 9).Value = 129310.125
$ws.Cells.Item(73, 10).Value = 4724.7
$ws.Cells.Item(73, 11).Value = 129310.125
$ws.Cells.Item(73, 12).Value = 4724.7
$ws.Cells.Item(73, 13).Value = -128374.125
$ws.Cells.Item(73, 14).Value = -6596.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 39303.152
$ws.Cells.Item(16, 9).Value = 46022.727
$ws.Cells.Item(16, 11).Value = 46022.727
$ws.Cells.Item(16, 13).Value = -45852.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3474.75
$ws.Cells.Item(46, 9).Value = 899
$ws.Cells.Item(46, 10).Value = 4333.3335
$ws.Cells.Item(46, 11).Value = 899
$ws.Cells.Item(46, 12).Value = 4333.3335
$ws.Cells.Item(46, 13).Value = -711
$ws.Cells.Item(46, 14).Value = -4709.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1571.4546
$ws.Cells.Item(82, 9).Value = 1428.6
$ws.Cells.Item(82, 10).Value = 3000
$ws.Cells.Item(82, 11).Value = 1428.6
$ws.Cells.Item(82, 12).Value = 3000
$ws.Cells.Item(82, 13).Value = -1067.6
$ws.Cells.Item(82, 14).Value = -3722

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1571.4546
$ws.Cells.Item(85, 9).Value = 1428.6
$ws.Cells.Item(85, 10).Value = 3000
$ws.Cells.Item(85, 11).Value = 1428.6
$ws.Cells.Item(85, 12).Value = 3000
$ws.Cells.Item(85, 13).Value = -180.5999999999999
$ws.Cells.Item(85, 14).Value = -5496

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 83655.75
$ws.Cells.Item(100, 9).Value = 111294.336
$ws.Cells.Item(100, 10).Value = 740
$ws.Cells.Item(100, 11).Value = 222588.672
$ws.Cells.Item(100, 12).Value = 1480
$ws.Cells.Item(100, 13).Value = -222047.672
$ws.Cells.Item(100, 14).Value = -2562
